# Update gait results after fixing downsampling issues
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update delta h (J column) raw input values with corrected downsampled data
$ws.Range("J2").Value = 0.068329323533756
$ws.Range("J3").Value = 0.0589235103059676
$ws.Range("J4").Value = 0.0806051206342671
$ws.Range("J5").Value = 0.0784181976776679
$ws.Range("J6").Value = 0.0764941143568182
$ws.Range("J7").Value = 0.0742879137475773
$ws.Range("J8").Value = 0.0925817852353351
$ws.Range("J9").Value = 0.0727276769579385
$ws.Range("J10").Value = 0.0467936513731346
$ws.Range("J11").Value = 0.0450165293024041
$ws.Range("J12").Value = 0.0505118583807846
$ws.Range("J13").Value = 0.0983252356080865
$ws.Range("J14").Value = 0.0796139458302568
$ws.Range("J15").Value = 0.086897943891136
$ws.Range("J16").Value = 0.0899088764756429
$ws.Range("J17").Value = 0.0997868288948745
$ws.Range("J18").Value = 0.0969537583739918
$ws.Range("J19").Value = 0.0593800862978098
$ws.Range("J20").Value = 0.0607622421758291
$ws.Range("J22").Value = 0.0668742051789103
$ws.Range("J23").Value = 0.0788429688679113
$ws.Range("J24").Value = 0.10018252919447
$ws.Range("J25").Value = 0.105745348348266
$ws.Range("J26").Value = 0.103106895737932
$ws.Range("J27").Value = 0.106513004311694
$ws.Range("J28").Value = 0.0776171479080649
$ws.Range("J29").Value = 0.064951671193058
$ws.Range("J30").Value = 0.0384968066453892
$ws.Range("J31").Value = 0.0241792109622463
$ws.Range("J32").Value = 0.0461085642434879
$ws.Range("J33").Value = 0.0800437105656894
$ws.Range("J34").Value = 0.0741156293980273
$ws.Range("J35").Value = 0.111369187408875
$ws.Range("J36").Value = 0.0855731847435321
$ws.Range("J37").Value = 0.103386352749795
$ws.Range("J38").Value = 0.0839255447245846
$ws.Range("J39").Value = 0.0734660291916763
$ws.Range("J41").Value = 0.0784227992708329
$ws.Range("J42").Value = 0.0903901020715342
$ws.Range("J43").Value = 0.0908105724372721
$ws.Range("J44").Value = 0.0984846713353026
$ws.Range("J45").Value = 0.10374508097948
$ws.Range("J46").Value = 0.0875992799139003
$ws.Range("J47").Value = 0.0951442122928678
$ws.Range("J48").Value = 0.0932320182061527
$ws.Range("J49").Value = 0.0672837176933147
$ws.Range("J50").Value = 0.0412309629559641
$ws.Range("J51").Value = 0.0403356009533752
$ws.Range("J52").Value = 0.0809216722409401
$ws.Range("J53").Value = 0.0752097923037162
$ws.Range("J54").Value = 0.0879647988792052
$ws.Range("J55").Value = 0.105212623742023
$ws.Range("J56").Value = 0.0773919149788697
$ws.Range("J57").Value = 0.0947042265554275
$ws.Range("J58").Value = 0.0669390872194007

# Recalculate all dependent formulas (AC:AH columns use J column values)
$excel.CalculateFull()

# Clear the bold/border/alignment formatting that was on the J1 header cell
$ws.Range("J1").ClearFormats()

# Leave the J column selected (whole-column selection), matching the saved view state
$ws.Columns("J:J").Select()
